# Add a new worksheet "ODI Batting Extra" after the existing "ODI Bowling" sheet,
# populate it with extra batting data (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH), and clean up 10 stray empty INNING_NUMBER
# cells left behind in the "ODI Batting" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Remove the stray empty B-column cells from "ODI Batting" (rows 7, 9, 11,
#    23, 25, 27, 28, 30, 32, 44) so those rows no longer carry a placeholder
#    empty INNING_NUMBER cell.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$emptyRows = @(7, 9, 11, 23, 25, 27, 28, 30, 32, 44)
foreach ($r in $emptyRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Batting Extra" sheet as the last tab.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

# Match the page margins used by the workbook's other sheets.
$extraSheet.PageSetup.LeftMargin = 54
$extraSheet.PageSetup.RightMargin = 54
$extraSheet.PageSetup.TopMargin = 72
$extraSheet.PageSetup.BottomMargin = 72
$extraSheet.PageSetup.HeaderMargin = 36
$extraSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 3) Header row: write the labels, then copy the bold/centered header style
#    from the "ODI Batting" sheet's header row onto the new sheet's header
#    row (must happen before the per-column NumberFormat tweak below, since
#    pasting formats resets column-level number formats).
# ---------------------------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extraSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$battingSheet.Range("A1:F1").Copy()
$extraSheet.Range("A1:F1").PasteSpecial(-4122)

# Columns A, C, D, E hold text-like values (zero-padded/percent strings) in
# the source data, so force text formatting on the data rows before writing
# them; column B holds a genuine number (batting position) and F is always
# the literal "NO" string (already non-numeric, no formatting needed).
$extraSheet.Range("A2:A46").NumberFormat = "@"
$extraSheet.Range("C2:C46").NumberFormat = "@"
$extraSheet.Range("D2:D46").NumberFormat = "@"
$extraSheet.Range("E2:E46").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 4) Data rows (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
#    PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH).
# ---------------------------------------------------------------------------
$rows = @(
    @(2, "3390", 6, "1", "3", "17.48%", "NO"),
    @(3, "3966", $null, $null, $null, $null, "NO"),
    @(4, "3967", 6, "0", "1", "4.58%", "NO"),
    @(5, "3968", $null, $null, $null, $null, "NO"),
    @(6, "3992", 8, "0", "3", "16.43%", "NO"),
    @(7, "3995", $null, $null, $null, $null, "NO"),
    @(8, "3997", $null, $null, $null, $null, "NO"),
    @(9, "3999", $null, $null, $null, $null, "NO"),
    @(10, "4002", $null, $null, $null, $null, "NO"),
    @(11, "4085", $null, $null, $null, $null, "NO"),
    @(12, "4088", $null, $null, $null, $null, "NO"),
    @(13, "4089", 7, "0", "0", "2.42%", "NO"),
    @(14, "4110", $null, $null, $null, $null, "NO"),
    @(15, "4114", 5, "2", "0", "10.70%", "NO"),
    @(16, "4137", $null, $null, $null, $null, "NO"),
    @(17, "4138", $null, $null, $null, $null, "NO"),
    @(18, "4139", $null, $null, $null, $null, "NO"),
    @(19, "4146", 6, "2", "2", "6.78%", "NO"),
    @(20, "4149", 7, "0", "0", "2.69%", "NO"),
    @(21, "4223", 7, "0", "0", $null, "NO"),
    @(22, "4225", $null, $null, $null, $null, "NO"),
    @(23, "4227", $null, $null, $null, $null, "NO"),
    @(24, "4242", $null, $null, $null, $null, "NO"),
    @(25, "4248", 8, $null, $null, $null, "NO"),
    @(26, "4249", 7, "2", "0", "5.07%", "NO"),
    @(27, "4250", 6, $null, $null, $null, "NO"),
    @(28, "4251", 6, $null, $null, $null, "NO"),
    @(29, "4252", 7, "4", "2", "11.21%", "NO"),
    @(30, "4305", 7, $null, $null, $null, "NO"),
    @(31, "4311", $null, $null, $null, $null, "NO"),
    @(32, "4315", 7, $null, $null, $null, "NO"),
    @(33, "4328", 7, "5", "2", "24.49%", "NO"),
    @(34, "4333", $null, $null, $null, $null, "NO"),
    @(35, "4337", $null, $null, $null, $null, "NO"),
    @(36, "4341", 6, "0", "0", $null, "NO"),
    @(37, "4346", $null, $null, $null, $null, "NO"),
    @(38, "4353", $null, $null, $null, $null, "NO"),
    @(39, "4355", $null, $null, $null, $null, "NO"),
    @(40, "4402", $null, $null, $null, $null, "NO"),
    @(41, "4406", $null, $null, $null, $null, "NO"),
    @(42, "4410", 7, "6", "3", "19.33%", "NO"),
    @(43, "4423", 7, "3", "0", "13.37%", "NO"),
    @(44, "4563", 7, $null, $null, $null, "NO"),
    @(45, "4566", 7, "1", "0", "6.06%", "NO"),
    @(46, "4568", $null, $null, $null, $null, "NO")
)

foreach ($row in $rows) {
    $r = $row[0]
    $extraSheet.Cells.Item($r, 1).Value = $row[1]
    if ($null -ne $row[2]) {
        $extraSheet.Cells.Item($r, 2).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $extraSheet.Cells.Item($r, 3).Value = $row[3]
    }
    if ($null -ne $row[4]) {
        $extraSheet.Cells.Item($r, 4).Value = $row[4]
    }
    if ($null -ne $row[5]) {
        $extraSheet.Cells.Item($r, 5).Value = $row[5]
    }
    $extraSheet.Cells.Item($r, 6).Value = $row[6]
}
